$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = Strike count) with recalculated values
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
